$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# 2018-12-17 as an Excel serial date number (matches the existing Issue Date column)
$issueDateSerial = 43451

# Row 13 (No=12)
$ws.Range("B13").Value = 17121812
$ws.Range("C13").Value = "Risk Management Plan.docx"
$ws.Range("D13").Value = "Balasubramanian"
$ws.Range("E2").Copy($ws.Range("E13"))
$ws.Range("E13").Value = $issueDateSerial

# Row 14 (No=13)
$ws.Range("B14").Value = 17121813
$ws.Range("C14").Value = "RAM_TURS.docx"
$ws.Range("D14").Value = "Balasubramanian"
$ws.Range("E2").Copy($ws.Range("E14"))
$ws.Range("E14").Value = $issueDateSerial

# Row 15 (No=14)
$ws.Range("B15").Value = 17121814
$ws.Range("C15").Value = "RAM_MRP.xlsx"
$ws.Range("D15").Value = "Balasubramanian"
$ws.Range("E2").Copy($ws.Range("E15"))
$ws.Range("E15").Value = $issueDateSerial

# Row 16 (No=15) - already formatted with border style; fill in values
$ws.Range("B16").Value = 17121815
$ws.Range("C16").Value = "STP.docx"
$ws.Range("D16").Value = "Nay Lin Aung"
$ws.Range("E2").Copy($ws.Range("E16"))
$ws.Range("E16").Value = $issueDateSerial

$ws.Range("J17").Select()
